# "2D array solved example"
#
# Row 7 of Sheet1 is the "Arrays" chapter row. This edit fills in the
# previously-empty "Question number" (G7) and "Comments" (H7) cells with
# details of the 2D-array solved example, un-minimizes the workbook
# window, and leaves the sheet scrolled/selected on the new comment.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Workbook window is no longer minimized.
$excel.WindowState = -4143   # xlNormal

# Fill in the question-number / comment cells for the 2D-array example.
$ws.Range("G7").Value = "10 problems"
$ws.Range("H7").Value = "6 from solved examples and 5 from exercise."

# Scroll so column D is at the left edge of the view, and leave the
# selection on the cell that was just edited (H7).
$win = $excel.ActiveWindow
$win.ScrollColumn = 4
$win.ScrollRow = 1
$ws.Range("H7").Select()
